# Atualizacao de bases das ligas, do dia: 21-04-2024 as 14:32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 118 and 119: the two match records swap places (B and F:AC) ---
# (A, C, D, E remain associated with the same row)

# Row 118 (id 116) gets new match data
$ws.Range('A118').Value = 116
$ws.Range('B118').Value = 7013885
$ws.Range('C118').Value = 'Uruguay Primera División'
$ws.Range('D118').Value = 'Uruguay Clausura'
$ws.Range('E118').Value = 45267.70833333334
$ws.Range('F118').Value = 'La Luz'
$ws.Range('G118').Value = 'Atletico Fenix Montevideo'
$ws.Range('H118').Value = 0
$ws.Range('I118').Value = 2
$ws.Range('J118').Value = 'A'
$ws.Range('K118').Value = 3
$ws.Range('L118').Value = 3
$ws.Range('M118').Value = 2.4
$ws.Range('N118').Value = 2.9
$ws.Range('O118').Value = 2.75
$ws.Range('P118').Value = 2.6
$ws.Range('Q118').Value = 0
$ws.Range('R118').Value = 2.025
$ws.Range('S118').Value = 1.825
$ws.Range('T118').Value = 2
$ws.Range('U118').Value = 2.025
$ws.Range('V118').Value = 1.825
$ws.Range('W118').Value = -1
$ws.Range('X118').Value = -1
$ws.Range('Y118').Value = 1.6
$ws.Range('Z118').Value = -1
$ws.Range('AA118').Value = 0.825
$ws.Range('AB118').Value = 0
$ws.Range('AC118').Value = -0

# Row 119 (id 117) gets new match data
$ws.Range('A119').Value = 117
$ws.Range('B119').Value = 7013886
$ws.Range('C119').Value = 'Uruguay Primera División'
$ws.Range('D119').Value = 'Uruguay Clausura'
$ws.Range('E119').Value = 45267.70833333334
$ws.Range('F119').Value = 'Racing Club de Montevideo'
$ws.Range('G119').Value = 'Cerro'
$ws.Range('H119').Value = 0
$ws.Range('I119').Value = 1
$ws.Range('J119').Value = 'A'
$ws.Range('K119').Value = 2.25
$ws.Range('L119').Value = 3.1
$ws.Range('M119').Value = 3.25
$ws.Range('N119').Value = 2.25
$ws.Range('O119').Value = 2.875
$ws.Range('P119').Value = 3.5
$ws.Range('Q119').Value = -0.25
$ws.Range('R119').Value = 1.95
$ws.Range('S119').Value = 1.9
$ws.Range('T119').Value = 2
$ws.Range('U119').Value = 1.925
$ws.Range('V119').Value = 1.925
$ws.Range('W119').Value = -1
$ws.Range('X119').Value = -1
$ws.Range('Y119').Value = 2.5
$ws.Range('Z119').Value = -1
$ws.Range('AA119').Value = 0.8999999999999999
$ws.Range('AB119').Value = -1
$ws.Range('AC119').Value = 0.925

# --- Rows 187 and 188: replace with new match data (add H, I, J stat columns) ---
# Row 187 (id 185) gets new match data
$ws.Range('A187').Value = 185
$ws.Range('B187').Value = 8081162
$ws.Range('C187').Value = 'Uruguay Primera División'
$ws.Range('D187').Value = 'Uruguay Apertura'
$ws.Range('E187').Value = 45402.41666666666
$ws.Range('F187').Value = 'Danubio'
$ws.Range('G187').Value = 'Cerro Largo'
$ws.Range('H187').Value = 1
$ws.Range('I187').Value = 1
$ws.Range('J187').Value = 'D'
$ws.Range('K187').Value = 2.3
$ws.Range('L187').Value = 3
$ws.Range('M187').Value = 3.4
$ws.Range('N187').Value = 2.4
$ws.Range('O187').Value = 2.9
$ws.Range('P187').Value = 3.25
$ws.Range('Q187').Value = -0.25
$ws.Range('R187').Value = 2.05
$ws.Range('S187').Value = 1.8
$ws.Range('T187').Value = 2
$ws.Range('U187').Value = 2
$ws.Range('V187').Value = 1.85
$ws.Range('W187').Value = -1
$ws.Range('X187').Value = 1.9
$ws.Range('Y187').Value = -1
$ws.Range('Z187').Value = -0.5
$ws.Range('AA187').Value = 0.4
$ws.Range('AB187').Value = 0
$ws.Range('AC187').Value = -0

# Row 188 (id 186) gets new match data
$ws.Range('A188').Value = 186
$ws.Range('B188').Value = 8081144
$ws.Range('C188').Value = 'Uruguay Primera División'
$ws.Range('D188').Value = 'Uruguay Apertura'
$ws.Range('E188').Value = 45402.625
$ws.Range('F188').Value = 'Boston River'
$ws.Range('G188').Value = 'Penarol'
$ws.Range('H188').Value = 1
$ws.Range('I188').Value = 3
$ws.Range('J188').Value = 'A'
$ws.Range('K188').Value = 4.75
$ws.Range('L188').Value = 3.75
$ws.Range('M188').Value = 1.727
$ws.Range('N188').Value = 4
$ws.Range('O188').Value = 3.6
$ws.Range('P188').Value = 1.909
$ws.Range('Q188').Value = 0.5
$ws.Range('R188').Value = 1.95
$ws.Range('S188').Value = 1.9
$ws.Range('T188').Value = 2.25
$ws.Range('U188').Value = 1.975
$ws.Range('V188').Value = 1.875
$ws.Range('W188').Value = -1
$ws.Range('X188').Value = -1
$ws.Range('Y188').Value = 0.909
$ws.Range('Z188').Value = -1
$ws.Range('AA188').Value = 0.8999999999999999
$ws.Range('AB188').Value = 0.9750000000000001
$ws.Range('AC188').Value = -1

# --- New rows 189 and 190: brand-new match records ---
# Copy cell formatting (border/bold style for id column, date format for date column)
$ws.Range('A188').Copy()
$ws.Range('A189:A190').PasteSpecial(-4122)
$ws.Range('E188').Copy()
$ws.Range('E189:E190').PasteSpecial(-4122)

# Row 189 (id 187) new match data
$ws.Range('A189').Value = 187
$ws.Range('B189').Value = 8081249
$ws.Range('C189').Value = 'Uruguay Primera División'
$ws.Range('D189').Value = 'Uruguay Apertura'
$ws.Range('E189').Value = 45402.75
$ws.Range('F189').Value = 'Nacional De Football'
$ws.Range('G189').Value = 'Rampla Juniors'
$ws.Range('H189').Value = 6
$ws.Range('I189').Value = 2
$ws.Range('J189').Value = 'H'
$ws.Range('K189').Value = 1.444
$ws.Range('L189').Value = 4
$ws.Range('M189').Value = 8.5
$ws.Range('N189').Value = 1.25
$ws.Range('O189').Value = 5
$ws.Range('P189').Value = 13
$ws.Range('Q189').Value = -1.5
$ws.Range('R189').Value = 1.85
$ws.Range('S189').Value = 2
$ws.Range('T189').Value = 2.5
$ws.Range('U189').Value = 1.925
$ws.Range('V189').Value = 1.925
$ws.Range('W189').Value = 0.25
$ws.Range('X189').Value = -1
$ws.Range('Y189').Value = -1
$ws.Range('Z189').Value = 0.8500000000000001
$ws.Range('AA189').Value = -1
$ws.Range('AB189').Value = 0.925
$ws.Range('AC189').Value = -1

# Row 190 (id 188) new match data
$ws.Range('A190').Value = 188
$ws.Range('B190').Value = 8081250
$ws.Range('C190').Value = 'Uruguay Primera División'
$ws.Range('D190').Value = 'Uruguay Apertura'
$ws.Range('E190').Value = 45402.85416666666
$ws.Range('F190').Value = 'Deportivo Maldonado'
$ws.Range('G190').Value = 'Miramar Misiones'
$ws.Range('H190').Value = 1
$ws.Range('I190').Value = 2
$ws.Range('J190').Value = 'A'
$ws.Range('K190').Value = 2.2
$ws.Range('L190').Value = 3.3
$ws.Range('M190').Value = 3.3
$ws.Range('N190').Value = 2.25
$ws.Range('O190').Value = 3.3
$ws.Range('P190').Value = 3.25
$ws.Range('Q190').Value = -0.25
$ws.Range('R190').Value = 1.975
$ws.Range('S190').Value = 1.875
$ws.Range('T190').Value = 2.25
$ws.Range('U190').Value = 1.9
$ws.Range('V190').Value = 1.95
$ws.Range('W190').Value = -1
$ws.Range('X190').Value = -1
$ws.Range('Y190').Value = 2.25
$ws.Range('Z190').Value = -1
$ws.Range('AA190').Value = 0.875
$ws.Range('AB190').Value = 0.8999999999999999
$ws.Range('AC190').Value = -1

# --- Rows 191 and 192: previously rows 187/188 data, shifted down with updated ids ---
$ws.Range('A188').Copy()
$ws.Range('A191:A192').PasteSpecial(-4122)
$ws.Range('E188').Copy()
$ws.Range('E191:E192').PasteSpecial(-4122)

# Row 191 (id 189) shifted-down former row 187 data
$ws.Range('A191').Value = 189
$ws.Range('B191').Value = 8081251
$ws.Range('C191').Value = 'Uruguay Primera División'
$ws.Range('D191').Value = 'Uruguay Apertura'
$ws.Range('E191').Value = 45403.54166666666
$ws.Range('F191').Value = 'Atletico Fenix Montevideo'
$ws.Range('G191').Value = 'Montevideo Wanderers'
$ws.Range('K191').Value = 2.5
$ws.Range('L191').Value = 3
$ws.Range('M191').Value = 3
$ws.Range('N191').Value = 2.5
$ws.Range('O191').Value = 2.9
$ws.Range('P191').Value = 3.1
$ws.Range('Q191').Value = -0.25
$ws.Range('R191').Value = 2.1
$ws.Range('S191').Value = 1.775
$ws.Range('T191').Value = 2
$ws.Range('U191').Value = 1.875
$ws.Range('V191').Value = 1.975
$ws.Range('W191').Value = 0
$ws.Range('X191').Value = 0
$ws.Range('Y191').Value = 0
$ws.Range('Z191').Value = 0
$ws.Range('AA191').Value = 0

# Row 192 (id 190) shifted-down former row 188 data
$ws.Range('A192').Value = 190
$ws.Range('B192').Value = 8081885
$ws.Range('C192').Value = 'Uruguay Primera División'
$ws.Range('D192').Value = 'Uruguay Apertura'
$ws.Range('E192').Value = 45403.64583333334
$ws.Range('F192').Value = 'CA River Plate'
$ws.Range('G192').Value = 'Club Atletico Progreso'
$ws.Range('K192').Value = 2.625
$ws.Range('L192').Value = 3.1
$ws.Range('M192').Value = 2.75
$ws.Range('N192').Value = 3.1
$ws.Range('O192').Value = 3.1
$ws.Range('P192').Value = 2.375
$ws.Range('Q192').Value = 0.25
$ws.Range('R192').Value = 1.8
$ws.Range('S192').Value = 2.05
$ws.Range('T192').Value = 2.25
$ws.Range('U192').Value = 1.975
$ws.Range('V192').Value = 1.875
$ws.Range('W192').Value = 0
$ws.Range('X192').Value = 0
$ws.Range('Y192').Value = 0
$ws.Range('Z192').Value = 0
$ws.Range('AA192').Value = 0

